$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 1 (Composants / Chef Module swap with Enseignant / Nombre d'heures)
$ws.Range("C1").Value = "Chef  Module"
$ws.Range("D1").Value = "Composants"

# Set new column widths (target stored OOXML widths: C=35, D=24.5703125)
$ws.Columns.Item(3).ColumnWidth = 34.16666666666666
$ws.Columns.Item(4).ColumnWidth = 23.736979166666668

# Update the active selection
$ws.Range("E8").Select()
